$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"

$ws.Range("C6").Select()
